$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: row 3 ("31478e65-...md") status (zh-cn/de-de columns) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: update Status cell (C3) text ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus

# Add new Error Detail (column L) value for row 3
$wsZh.Range("L3").Value = "Handback file name: larbvcgz.prm is different with handoff file name: 31478e65-16f3-41f4-8eb8-0283831b774b.63d1cfb9b52f5fbef610af74dad4d2d59fbc6c91.zh-cn."

# --- de-de sheet: update Status cell (C3) text and add new Error Detail (column L) value for row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("L3").Value = "Handback file name: larbvcgz.prm is different with handoff file name: 31478e65-16f3-41f4-8eb8-0283831b774b.63d1cfb9b52f5fbef610af74dad4d2d59fbc6c91.de-de."
